$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the "dump_acc_O3 / every_dump / every_dump2" block (D:F) into a
# new G:I block (same formatting/style), giving us a second set of columns
# to hold the "Buffered" benchmark numbers.
$ws.Range("D2:F7").Copy($ws.Range("G2"))

# Overwrite the freshly-copied data cells (rows 3-7) with the new "Buffered"
# benchmark timings.
$ws.Range("G3").Value = "0.670s"
$ws.Range("H3").Value = "1.420s"
$ws.Range("I3").Value = "1.452s"

$ws.Range("G4").Value = "10.290s"
$ws.Range("H4").Value = "11.481s"
$ws.Range("I4").Value = "11.461s"

$ws.Range("G5").Value = "0.557s"
$ws.Range("H5").Value = "1.477s"
$ws.Range("I5").Value = "1.436s"

$ws.Range("G6").Value = "0.500s"
$ws.Range("H6").Value = "1.472s"
$ws.Range("I6").Value = "1.445s"

$ws.Range("G7").Value = "0.488s"
$ws.Range("H7").Value = "1.467s"
$ws.Range("I7").Value = "1.428s"

# Add a new super-header row above the column headers labelling the two
# blocks as "Unbuffered" (original D:F columns) and "Buffered" (new G:I
# columns).
$ws.Range("D1").Value = "Unbuffered"
$ws.Range("G1").Value = "Buffered"

$ws.Range("D1").Font.Name = "Arial"
$ws.Range("D1").Font.ThemeColor = 1

$ws.Range("G1").Font.Name = "Arial"
$ws.Range("G1").Font.ThemeColor = 1
